# Actualización automática del mapa (2025-09-04 07:53:20)
# Agrega la nueva fila de datos (caso -581) al final de la hoja "AYKO".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

# Las columnas A-L, O, P son texto (incluso cuando el contenido parece
# numérico, p.ej. "-581"), por lo que forzamos el formato de texto antes
# de escribir el valor para que Excel no las reinterprete como números o
# fechas.
$textCols = @(1,2,3,5,6,7,8,10,11,12,15,16)
foreach ($col in $textCols) {
    $ws.Cells.Item($newRow, $col).NumberFormat = "@"
}

$ws.Cells.Item($newRow, 1).Value  = "-581"
$ws.Cells.Item($newRow, 2).Value  = "9/4/2025"
$ws.Cells.Item($newRow, 3).Value  = "Praga 1380"
# Columna D (Comuna) queda vacía para este registro.
$ws.Cells.Item($newRow, 5).Value  = "809432814"
$ws.Cells.Item($newRow, 6).Value  = "AYKO"
$ws.Cells.Item($newRow, 7).Value  = "Pendiente"
$ws.Cells.Item($newRow, 8).Value  = "Cambiar"
$ws.Cells.Item($newRow, 9).Value  = 0
$ws.Cells.Item($newRow, 10).Value = "Cambio"
$ws.Cells.Item($newRow, 11).Value = "Sin equipos"
$ws.Cells.Item($newRow, 12).Value = "Pasante"
$ws.Cells.Item($newRow, 13).Value = -58.481608
$ws.Cells.Item($newRow, 14).Value = -34.587837
$ws.Cells.Item($newRow, 15).Value = "Paternal"
$ws.Cells.Item($newRow, 16).Value = "Capital Norte"
